$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price & 1h volume change) per latest scrape
$ws.Range("D2").Value = "42.011.63"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.220.84"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.85"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.14"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.617"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.10"
$ws.Range("E10").Value = "  +4.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0963"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.16"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "2.554.40"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.30"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.844"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "2.215.39"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "41.898.55"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("E19").Value = "  +11.93%  "
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.55"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.88"
$ws.Range("E22").Value = "  +31.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.50"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("E24").Value = "  -7.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.60"
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.84"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("E32").Value = "  +15.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0803"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.42"
$ws.Range("E35").Value = "  -3.72%  "
$ws.Range("E36").Value = "  -4.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.31"
$ws.Range("E37").Value = "  -5.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0303"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.07"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.14"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.37"
$ws.Range("E41").Value = "  +6.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.66"
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.201"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.38"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  +6.89%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.428.54"
$ws.Range("E51").Value = "  -1.31%  "
